$wb = $excel.ActiveWorkbook
$ws2017 = $wb.Worksheets.Item("2017")

# --- 1) Create the "2018" sheet as a copy of "2017" ---
$ws2017.Copy([System.Reflection.Missing]::Value, $ws2017)
$ws2018 = $wb.Worksheets.Item("2017 (2)")
$ws2018.Name = "2018"

# --- 2) On "2018": G column now pulls last year's (2017) full-year total (U col) ---
for ($row = 2; $row -le 18; $row++) {
    $ws2018.Range("G$row").Formula = "='2017'!U$row"
}

# --- 3) On "2018": zero out all monthly dividend entries (H:S) for the new year ---
$ws2018.Range("H2:S18").Value2 = 0
# one dividend already posted early in the new year
$ws2018.Range("H17").Value2 = 7.92

# --- 4) On "2018": sheet view / selection ---
$ws2018.Activate()
$ws2018.Range("D18").Select()

# --- 5) On "2017": convert the F column (buy totals from 'Initial Buys') from
#        formulas to plain cached values ---
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws2017.Range("F$row")
    $cell.Value2 = $cell.Value2
}
# a few of those pasted values were hand-rounded by the author
$ws2017.Range("F6").Value2 = 1625.82
$ws2017.Range("F13").Value2 = 1173.32
$ws2017.Range("F16").Value2 = 1356.18

# --- 6) On "2017": real data edits reflected in the commit ---
$ws2017.Range("D17").Value2 = 15.523999999999999
$ws2017.Range("S17").Value2 = 0

# --- 7) On "2017": sheet view / selection ---
$ws2017.Activate()
$ws2017.Range("F3").Select()

# --- 8) Leave "2018" as the active tab (it is the most-recently added sheet) ---
$ws2018.Activate()

Write-Output "done"
